$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A in this sheet currently has 7 vertically-merged blocks (one per
# binding-site group). The edit un-merges each block and repeats the group
# label down every row, changing the alignment from centered (both axes) to
# left/general + vertically centered, while keeping the bold font.
$ranges = @("A2:A33", "A34:A77", "A78:A95", "A96:A105", "A106:A116", "A117:A127", "A128:A129")

foreach ($rangeAddr in $ranges) {
    $rng = $ws.Range($rangeAddr)
    $label = $rng.Cells.Item(1, 1).Value2
    $rng.UnMerge()
    $rng.Value2 = $label
    $rng.Font.Bold = $true
    $rng.HorizontalAlignment = 1
    $rng.VerticalAlignment = -4108
}

# Restore the selection left behind by the author's last interactive edit.
$ws.Range("G11").Select() | Out-Null
